$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "28.000.07"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.859.81"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.21%  "
Set-TextValue $ws.Range("D5") "311.80"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.20%  "
Set-TextValue $ws.Range("D7") "0.5118"
$ws.Range("E7").Value = "  +2.41%  "
Set-TextValue $ws.Range("D8") "0.3808"
$ws.Range("E8").Value = "  -0.97%  "
Set-TextValue $ws.Range("D9") "0.08314"
$ws.Range("E9").Value = "  -9.37%  "
Set-TextValue $ws.Range("D10") "1.106"
$ws.Range("E10").Value = "  -1.18%  "
Set-TextValue $ws.Range("D11") "41.29"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D13") "20.43"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.857.89"
$ws.Range("E14").Value = "  -0.93%  "
Set-TextValue $ws.Range("D15") "7.178"
$ws.Range("E15").Value = "  -1.44%  "
Set-TextValue $ws.Range("D16") "1.003"
$ws.Range("E16").Value = "  +0.16%  "
Set-TextValue $ws.Range("D17") "0.00001094"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("E18").Value = "  -1.13%  "
Set-TextValue $ws.Range("D19") "0.06615"
$ws.Range("E19").Value = "  -0.24%  "
Set-TextValue $ws.Range("D20") "17.81"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +0.15%  "
Set-TextValue $ws.Range("D22") "6.007"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").Value = "28.028.90"
$ws.Range("E23").Value = "  -0.03%  "
Set-TextValue $ws.Range("D24") "11.04"
$ws.Range("E24").Value = "  -2.93%  "
Set-TextValue $ws.Range("D25") "2.254"
$ws.Range("E25").Value = "  -2.36%  "
Set-TextValue $ws.Range("D26") "2.568"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").Value = "2.076.98"
$ws.Range("E27").Value = "  -0.74%  "
Set-TextValue $ws.Range("D28") "157.27"
$ws.Range("E28").Value = "  -0.22%  "
Set-TextValue $ws.Range("D29") "20.48"
$ws.Range("E29").Value = "  -1.33%  "
Set-TextValue $ws.Range("D30") "124.58"
$ws.Range("E30").Value = "  -1.52%  "
Set-TextValue $ws.Range("D31") "0.1059"
$ws.Range("E31").Value = "  +0.43%  "
Set-TextValue $ws.Range("D32") "1.040"
$ws.Range("E32").Value = "  -2.48%  "
Set-TextValue $ws.Range("D33") "5.588"
$ws.Range("E33").Value = "  +0.10%  "
Set-TextValue $ws.Range("D34") "3.605"
$ws.Range("E34").Value = "  +0.43%  "
Set-TextValue $ws.Range("D35") "9.604"
$ws.Range("E35").Value = "  +2.73%  "
Set-TextValue $ws.Range("D36") "0.06516"
$ws.Range("E36").Value = "  -0.77%  "
Set-TextValue $ws.Range("D37") "0.02410"
$ws.Range("E37").Value = "  +0.48%  "
Set-TextValue $ws.Range("D38") "0.2153"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  +0.38%  "
Set-TextValue $ws.Range("D40") "0.6399"
$ws.Range("E40").Value = "  +0.04%  "
Set-TextValue $ws.Range("D41") "1.231"
$ws.Range("E41").Value = "  -4.15%  "
Set-TextValue $ws.Range("D42") "11.25"
$ws.Range("E42").Value = "  -3.12%  "
Set-TextValue $ws.Range("D43") "4.866"
$ws.Range("E43").Value = "  -1.31%  "
Set-TextValue $ws.Range("D44") "0.6061"
$ws.Range("E44").Value = "  +0.86%  "
Set-TextValue $ws.Range("D45") "13.09"
$ws.Range("E45").Value = "  -1.66%  "
Set-TextValue $ws.Range("D46") "1.283"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -0.75%  "
Set-TextValue $ws.Range("D49") "1.207"
$ws.Range("E49").Value = "  +0.45%  "
Set-TextValue $ws.Range("D50") "120.42"
$ws.Range("E50").Value = "  -0.67%  "
Set-TextValue $ws.Range("D51") "79.60"
